# Natmi following Dr Hou advice
# The clustering gained a new "ECs" cluster, which expands the Wnt1-Fzd9
# LR-pair table from 2 data rows to 6 data rows (rows 2-7): every
# combination of sending cluster (ECs, FAPs) x target cluster (ECs, FAPs,
# sCs), each with recomputed NATMI expression/specificity statistics.
# Ligand/receptor symbols (columns B, C) remain Wnt1 / Fzd9 throughout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, then values for columns A..T (20 columns).
$rows = @(
    , @(2, @("ECs", "Wnt1", "Fzd9", "ECs",
             1, 0.3333333333333333, 0.027123, 0.081369,
             0.07131444737854614, 0.07131444737854616, 2, 0.6666666666666666,
             0.4056473333333333, 1.216942, 0.2120329281749088, 0.2120329281749088,
             0.011002372622, 0.09902135359799999, 0.01512101109884859, 0.01512101109884859))
    , @(3, @("ECs", "Wnt1", "Fzd9", "FAPs",
             1, 0.3333333333333333, 0.027123, 0.081369,
             0.07131444737854614, 0.07131444737854616, 3, 1,
             1.086432, 3.259296, 0.5678808642225905, 0.5678808642225904,
             0.029467295136, 0.265205656224, 0.04049811000888524, 0.04049811000888524))
    , @(4, @("ECs", "Wnt1", "Fzd9", "sCs",
             1, 0.3333333333333333, 0.027123, 0.081369,
             0.07131444737854614, 0.07131444737854616, 3, 1,
             0.4210543333333334, 1.263163, 0.2200862076025007, 0.2200862076025006,
             0.011420256683, 0.102782310147, 0.01569532627081232, 0.01569532627081232))
    , @(5, @("FAPs", "Wnt1", "Fzd9", "ECs",
             3, 1, 0.3532066666666667, 1.05962,
             0.9286855526214538, 0.9286855526214538, 2, 0.6666666666666666,
             0.4056473333333333, 1.216942, 0.2120329281749088, 0.2120329281749088,
             0.1432773424488889, 1.28949608204, 0.1969119170760602, 0.1969119170760602))
    , @(6, @("FAPs", "Wnt1", "Fzd9", "FAPs",
             3, 1, 0.3532066666666667, 1.05962,
             0.9286855526214538, 0.9286855526214538, 3, 1,
             1.086432, 3.259296, 0.5678808642225905, 0.5678808642225904,
             0.38373502528, 3.45361522752, 0.5273827542137053, 0.5273827542137052))
    , @(7, @("FAPs", "Wnt1", "Fzd9", "sCs",
             3, 1, 0.3532066666666667, 1.05962,
             0.9286855526214538, 0.9286855526214538, 3, 1,
             0.4210543333333334, 1.263163, 0.2200862076025007, 0.2200862076025006,
             0.1487191975622222, 1.33847277806, 0.2043908813316884, 0.2043908813316883))
)

foreach ($entry in $rows) {
    $r = $entry[0]
    $vals = $entry[1]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $vals[$c]
    }
}
